# WebForm User Assignment execution
# Updates the generated phone-number values (column F) and the
# Match1UserPos value (AM2) on Sheet1 to reflect a fresh run of the
# WebForm automation (TC028_WF_RRY_OneYN_TwoN_Test).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (shared string), matching the sheet's
# existing convention for this column, instead of letting Excel
# auto-convert the numeric-looking phone numbers into numbers. A leading
# apostrophe forces text entry (standard Excel behavior); the style is
# then reset back to the sheet's default (same as neighboring cells) so
# no stray formatting is introduced.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = $ws.Range("A1").Style
}

# Column F (PN_Value) - phone number used for this test run
Set-TextValue "F2"  "9840000680"
Set-TextValue "F3"  "9840049855"
Set-TextValue "F4"  "9840034871"
Set-TextValue "F5"  "9840037996"
Set-TextValue "F6"  "9840097183"
Set-TextValue "F7"  "9840096257"
Set-TextValue "F8"  "9840098514"
Set-TextValue "F9"  "9840081325"
Set-TextValue "F10" "9840087970"

# Match1UserPos result for row 2
Set-TextValue "AM2" "0"
